# Fruta / hortaliza, semanal
#
# A new weekly price-report row needs to be inserted into the "Tuna"
# sheet right above the existing row 19 (pushing rows 19-48 down to
# 20-49) and filled in with this week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 19; this shifts the current rows
# 19..48 down to 20..49 and keeps all of their data/formatting intact.
$ws.Rows.Item(19).Insert()

# Fill in the newly inserted row 19 with the new weekly record.
$ws.Range("A19").Value = 5
$ws.Range("B19").Value = "Macroferia Regional de Talca"
$ws.Range("C19").Value = "Maule"
$ws.Range("D19").Value = 45177
$ws.Range("E19").Value = 7
$ws.Range("F19").Value = "Fruta"
$ws.Range("G19").Value = 100107
$ws.Range("H19").Value = "Otros"
$ws.Range("I19").Value = 100107011
$ws.Range("J19").Value = "Tuna"
$ws.Range("K19").Value = "Sin especificar"
$ws.Range("L19").Value = "Primera"
$ws.Range("M19").Value = 50
$ws.Range("N19").Value = 24000
$ws.Range("O19").Value = 24000
$ws.Range("P19").Value = 24000
$ws.Range("Q19").Value = "`$/caja 18 kilos"
$ws.Range("R19").Value = "Provincia de Limarí"
$ws.Range("S19").Value = 1333
$ws.Range("T19").Value = 18
